$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Turn score): Mean updated
$ws.Range("D3").Value = 124.60358675391

# Row 4 (Turn nb roll): Max and Mean updated
$ws.Range("C4").Value = 5850
$ws.Range("D4").Value = 435.97225

# Row 5 (Turn nb full roll): Mean updated
$ws.Range("D5").Value = 3.498874

# Row 6 (Turn nb bonus): Max and Mean updated
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 0.4497412596166652

# Row 7 (Roll nb dice fail roll): Mean updated
$ws.Range("D7").Value = 0.338133

# Row 8 (Roll nb dice to roll): Mean updated
$ws.Range("D8").Value = 2.173811

# Row 9 (-): Mean updated
$ws.Range("D9").Value = 3.34426818456452
